# The two worksheets swap roles/content:
#   - sheet #1 (physically xl/worksheets/sheet1.xml, currently named
#     "hotel_info" and holding the hotel_info header+data) becomes the
#     "review_info" sheet (header row only, 25 columns).
#   - sheet #2 (physically xl/worksheets/sheet2.xml, currently named
#     "review_info" and holding only a header row) becomes the
#     "hotel_info" sheet (header + one data row, 10 columns — a new
#     "State" column is inserted right after "Hotel_Name").

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item(1)   # currently "hotel_info"
$wsReview = $wb.Worksheets.Item(2)   # currently "review_info"

# Copy the three numeric-looking text values (English_Reviews_num,
# Local_Rank, Total_Reviews_num = "459"/"14"/"467") from their old
# G2:I2 spot on the hotel sheet straight across to their new H2:J2
# spot on the (soon to be) hotel_info sheet FIRST, via copy/paste, so
# they keep their original text cell-type instead of being
# re-interpreted as numbers by a plain value assignment.
$wsHotel.Range("G2:I2").Copy()
$wsReview.Range("H2:J2").PasteSpecial()

# Now it is safe to wipe both sheets' remaining old content (except
# the H2:J2 values we just placed on $wsReview).
$wsHotel.Cells.Clear()
$wsReview.Range("A1:G1").Clear()
$wsReview.Range("K1:Y1").Clear()

# ---- header row for the new "review_info" sheet (25 columns) ----
$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)
for ($c = 1; $c -le $reviewHeaders.Length; $c++) {
    $wsHotel.Cells.Item(1, $c).Value = $reviewHeaders[$c - 1]
}

# ---- header row for the new "hotel_info" sheet (10 columns) ----
$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)
for ($c = 1; $c -le $hotelHeaders.Length; $c++) {
    $wsReview.Cells.Item(1, $c).Value = $hotelHeaders[$c - 1]
}

# ---- data row for the new "hotel_info" sheet ----
$wsReview.Cells.Item(2, 1).Value = 5148
$wsReview.Cells.Item(2, 2).Value = "La Quinta Inns & Suites New Orleans Veterans Metarie"
$wsReview.Cells.Item(2, 3).Value = "Louisiana"
$wsReview.Cells.Item(2, 4).Value = "Metairie"
$wsReview.Cells.Item(2, 5).Value = 70003
$wsReview.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40314-d93084-Reviews-La_Quinta_Inn_New_Orleans_Veterans_Metairie-Metairie_Louisiana.html"
$wsReview.Cells.Item(2, 7).Value = "La Quinta Inn New Orleans Veterans / Metairie"
# columns 8,9,10 (H2:J2 = "459"/"14"/"467") were already populated above via copy/paste

# ---- rename the sheets to match their new contents ----
# Route through a throwaway intermediate name so the two renames never
# collide (both sheets can't share a name even for an instant).
$wsHotel.Name = "swap_tmp_name"
$wsReview.Name = "hotel_info"
$wsHotel.Name = "review_info"
